$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("test", "page1", "Страница 1", "Page one"),
    @("test", "page2", "Страница 2", "Page two"),
    @("nv", "home", "Главная", "Start page"),
    @("nv", "participate", "Участвовать!", "Participate!"),
    @("nv", "statistics", "Статистика", "Statistics"),
    @("nv", "about", "О проекте", "About"),
    @("nv", "team", "Наша команда", "Our team"),
    @("nv", "for_volunteers", "Для волонтеров", "For volunteers"),
    @("nv", "profit_science", "Польза для науки", "Profit for science"),
    @("nv", "profit_personal", "Польза для вас", "Profit for you"),
    @("nv", "howtohelp", "Как нам помочь", "How to help"),
    @("nv", "voluntary_project", "Наш волонтерский проект", "Our voluntary project"),
    @("nv", "for_scientists", "Для специалистов", "For scientists"),
    @("nv", "cooperation", "Сотрудничество", "Cooperation"),
    @("nv", "web_app", "Наше веб-приложение", "Our web-application"),
    @("nv", "scientific_project", "Наш научный проект", "Our scientific project"),
)

$startRow = 3
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowVals[$c]
    }
}

$ws.Range("B17").Select()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 4

Write-Host "Done filling navbar translations"
